# Natmi following Dr Hou advice
# Update the LR-pairs_lrc2p Wnt2b-Fzd4 sheet with recomputed NATMI statistics
# for rows 2-16 (data rows). Columns F (Ligand detection rate) and L
# (Receptor detection rate) are unchanged (always 1); all other numeric
# columns from E through T are refreshed with the newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices (1-based) for: E, G, H, I, J, K, M, N, O, P, Q, R, S, T
$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)
$data = @(
    ,@(3, 1.584191, 4.752573, 0.2039202590281707, 0.2147775532998296, 3, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 46.30336083150867, 416.730247483578, 0.09078777319161138, 0.1071423406249714)
    ,@(3, 1.584191, 4.752573, 0.2039202590281707, 0.2147775532998296, 3, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 24.14997685663267, 217.349791709694, 0.04735126310638461, 0.05588114987752189)
    ,@(3, 1.584191, 4.752573, 0.2039202590281707, 0.2147775532998296, 2, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 33.549580355767, 201.297482134602, 0.06578122273017462, 0.05175406279733635)
    ,@(3, 2.213251666666667, 6.639755, 0.2848942161400975, 0.3000627940718238, 3, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 64.68979468549223, 582.20815216943, 0.1268383612388211, 0.1496871046223502)
    ,@(3, 2.213251666666667, 6.639755, 0.2848942161400975, 0.3000627940718238, 3, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 33.73960370176556, 303.6564333158901, 0.06615380467947211, 0.0780707932955528)
    ,@(3, 2.213251666666667, 6.639755, 0.2848942161400975, 0.3000627940718238, 2, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 46.87166171147834, 281.22997026887, 0.0919020502218042, 0.07230489615392083)
    ,@(3, 1.261258666666667, 3.783776, 0.1623517580949468, 0.1709958272107795, 3, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 36.86456692692623, 331.7811023423361, 0.07228097228508906, 0.08530171278601963)
    ,@(3, 1.261258666666667, 3.783776, 0.1623517580949468, 0.1709958272107795, 3, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 19.22708032694756, 173.043722942528, 0.03769885763177621, 0.04448995391737702)
    ,@(3, 1.261258666666667, 3.783776, 0.1623517580949468, 0.1709958272107795, 2, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 26.71060433163733, 160.263625989824, 0.05237192817808148, 0.04120416050738288)
    ,@(3, 1.531826, 4.595478, 0.1971797306676319, 0.207678140048179, 3, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 44.77281590987867, 402.955343188908, 0.08778680819232863, 0.1036007798745147)
    ,@(3, 1.531826, 4.595478, 0.1971797306676319, 0.207678140048179, 3, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 23.35170598014267, 210.165353821284, 0.04578608006181118, 0.05403401375988428)
    ,@(3, 1.531826, 4.595478, 0.1971797306676319, 0.207678140048179, 2, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 32.440608157762, 194.643648946572, 0.0636068424134921, 0.05004334641377999)
    ,@(2, 1.1781515, 2.356303, 0.1516540360691531, 0.106485685369388, 3, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 34.43547780455967, 206.612866827358, 0.06751821665907504, 0.053120660880252)
    ,@(2, 1.1781515, 2.356303, 0.1516540360691531, 0.106485685369388, 3, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 17.96016481510567, 107.760988890634, 0.03521479522082987, 0.02770560727838467)
    ,@(2, 1.1781515, 2.356303, 0.1516540360691531, 0.106485685369388, 2, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 24.9505826131555, 99.802330452622, 0.04892102418924821, 0.02565941721075132)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}
